# Auto-generated edit script: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) for specific Leve rows across all 8 job sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 366.125
$ws.Range("J2").Value = 447.5
$ws.Range("L2").Value = 447.5
$ws.Range("N2").Value = -673.5
# Row 28
$ws.Range("H28").Value = 53886.79
$ws.Range("I28").Value = 67360.60000000001
$ws.Range("K28").Value = 67360.60000000001
$ws.Range("M28").Value = -66875.60000000001
# Row 52
$ws.Range("H52").Value = 1492.8572
$ws.Range("I52").Value = 816.6667
$ws.Range("K52").Value = 2450.0001
$ws.Range("M52").Value = -2290.0001
# Row 86
$ws.Range("H86").Value = 11126538
$ws.Range("I86").Value = 9140.75
$ws.Range("J86").Value = 33361334
$ws.Range("K86").Value = 9140.75
$ws.Range("L86").Value = 33361334
$ws.Range("M86").Value = -8017.75
$ws.Range("N86").Value = -33363580
# Row 89
$ws.Range("H89").Value = 11126538
$ws.Range("I89").Value = 9140.75
$ws.Range("J89").Value = 33361334
$ws.Range("K89").Value = 45703.75
$ws.Range("L89").Value = 166806670
$ws.Range("M89").Value = -40087.75
$ws.Range("N89").Value = -166817902
# Row 129
$ws.Range("H129").Value = 1243.0714
$ws.Range("I129").Value = 785.8
$ws.Range("K129").Value = 2357.4
$ws.Range("M129").Value = 2642.6
# Row 138
$ws.Range("H138").Value = 2960.4792
$ws.Range("I138").Value = 1585.6923
$ws.Range("J138").Value = 3471.1143
$ws.Range("K138").Value = 4757.0769
$ws.Range("L138").Value = 10413.3429
$ws.Range("M138").Value = 382.9231
$ws.Range("N138").Value = -20693.3429
# Row 141
$ws.Range("H141").Value = 2935.8635
$ws.Range("I141").Value = 2635.7058
$ws.Range("K141").Value = 7907.117400000001
$ws.Range("M141").Value = -2727.117400000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1283.9701
$ws.Range("I32").Value = 1283.9701
$ws.Range("K32").Value = 1283.9701
$ws.Range("M32").Value = -996.9701
# Row 45
$ws.Range("H45").Value = 1753.6471
$ws.Range("I45").Value = 1658
$ws.Range("K45").Value = 1658
$ws.Range("M45").Value = -1281
# Row 74
$ws.Range("H74").Value = 1684.7222
$ws.Range("I74").Value = 1023
$ws.Range("K74").Value = 1023
$ws.Range("M74").Value = -149
# Row 77
$ws.Range("H77").Value = 1684.7222
$ws.Range("I77").Value = 1023
$ws.Range("K77").Value = 5115
$ws.Range("M77").Value = -747
# Row 110
$ws.Range("H110").Value = 66734930
$ws.Range("I110").Value = 90955820
$ws.Range("K110").Value = 90955820
$ws.Range("M110").Value = -90953775
# Row 122
$ws.Range("H122").Value = 16669061
$ws.Range("I122").Value = 17546116
$ws.Range("K122").Value = 52638348
$ws.Range("M122").Value = -52635898
# Row 132
$ws.Range("H132").Value = 43480760
$ws.Range("I132").Value = 62501984
$ws.Range("K132").Value = 187505952
$ws.Range("M132").Value = -187503422

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 18
$ws.Range("H18").Value = 2981.3333
$ws.Range("J18").Value = 2981.3333
$ws.Range("L18").Value = 2981.3333
$ws.Range("N18").Value = -4039.3333
# Row 20
$ws.Range("H20").Value = 12594.44
$ws.Range("I20").Value = 14948.333
$ws.Range("K20").Value = 14948.333
$ws.Range("M20").Value = -14701.333
# Row 94
$ws.Range("H94").Value = 3205.1924
$ws.Range("I94").Value = 2878.8635
$ws.Range("K94").Value = 2878.8635
$ws.Range("M94").Value = -2427.8635
# Row 99
$ws.Range("H99").Value = 884.5
$ws.Range("I99").Value = 786.8333
$ws.Range("K99").Value = 786.8333
$ws.Range("M99").Value = 711.1667
# Row 134
$ws.Range("H134").Value = 2621.9062
$ws.Range("I134").Value = 2539.25
$ws.Range("J134").Value = 3200.5
$ws.Range("K134").Value = 7617.75
$ws.Range("L134").Value = 9601.5
$ws.Range("M134").Value = -5082.75
$ws.Range("N134").Value = -14671.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 86.84614999999999
$ws.Range("I7").Value = 93.25
$ws.Range("J7").Value = 76.59999999999999
$ws.Range("K7").Value = 93.25
$ws.Range("L7").Value = 76.59999999999999
$ws.Range("M7").Value = 19.75
$ws.Range("N7").Value = -302.6
# Row 28
$ws.Range("H28").Value = 35700
$ws.Range("J28").Value = 35700
$ws.Range("L28").Value = 35700
$ws.Range("N28").Value = -36190
# Row 99
$ws.Range("H99").Value = 2644.4443
$ws.Range("I99").Value = 2300
$ws.Range("K99").Value = 2300
$ws.Range("M99").Value = -802
# Row 126
$ws.Range("H126").Value = 2644.4443
$ws.Range("I126").Value = 2300
$ws.Range("K126").Value = 6900
$ws.Range("M126").Value = -4430
# Row 132
$ws.Range("H132").Value = 2812.2222
$ws.Range("I132").Value = 2606.682
$ws.Range("K132").Value = 7820.045999999999
$ws.Range("M132").Value = -5290.045999999999
# Row 140
$ws.Range("H140").Value = 59999
$ws.Range("J140").Value = 59999
$ws.Range("L140").Value = 59999
$ws.Range("N140").Value = -70359
# Row 141
$ws.Range("H141").Value = 108347.6
$ws.Range("J141").Value = 108347.6
$ws.Range("L141").Value = 108347.6
$ws.Range("N141").Value = -118707.6

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 940.75
$ws.Range("I55").Value = 143.92857
$ws.Range("J55").Value = 2800
$ws.Range("K55").Value = 431.78571
$ws.Range("L55").Value = 8400
$ws.Range("M55").Value = -254.78571
$ws.Range("N55").Value = -8754
# Row 70
$ws.Range("H70").Value = 1500
$ws.Range("I70").Value = 1500
$ws.Range("K70").Value = 4500
$ws.Range("M70").Value = -4185
# Row 73
$ws.Range("H73").Value = 1500
$ws.Range("I73").Value = 1500
$ws.Range("K73").Value = 4500
$ws.Range("M73").Value = -3408

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 17426.625
$ws.Range("J24").Value = 17426.625
$ws.Range("L24").Value = 17426.625
$ws.Range("N24").Value = -17772.625
# Row 80
$ws.Range("H80").Value = 3304.7307
$ws.Range("I80").Value = 3419.875
$ws.Range("J80").Value = 3253.5557
$ws.Range("K80").Value = 3419.875
$ws.Range("L80").Value = 3253.5557
$ws.Range("M80").Value = -2421.875
$ws.Range("N80").Value = -5249.5557
# Row 83
$ws.Range("H83").Value = 3304.7307
$ws.Range("I83").Value = 3419.875
$ws.Range("J83").Value = 3253.5557
$ws.Range("K83").Value = 17099.375
$ws.Range("L83").Value = 16267.7785
$ws.Range("M83").Value = -12107.375
$ws.Range("N83").Value = -26251.7785
# Row 126
$ws.Range("H126").Value = 9806
$ws.Range("J126").Value = 3670.7144
$ws.Range("L126").Value = 11012.1432
$ws.Range("N126").Value = -15952.1432
# Row 132
$ws.Range("H132").Value = 2857.8718
$ws.Range("I132").Value = 1866.6923
$ws.Range("K132").Value = 5600.0769
$ws.Range("M132").Value = -3070.0769

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 950000
$ws.Range("J4").Value = 900000
$ws.Range("L4").Value = 900000
$ws.Range("N4").Value = -900226
# Row 28
$ws.Range("H28").Value = 950000
$ws.Range("J28").Value = 900000
$ws.Range("L28").Value = 900000
$ws.Range("N28").Value = -900464
# Row 37
$ws.Range("H37").Value = 950000
$ws.Range("J37").Value = 900000
$ws.Range("L37").Value = 900000
$ws.Range("N37").Value = -900214
# Row 132
$ws.Range("H132").Value = 3297.8
$ws.Range("I132").Value = 2988.5
$ws.Range("J132").Value = 3691.4546
$ws.Range("K132").Value = 8965.5
$ws.Range("L132").Value = 11074.3638
$ws.Range("M132").Value = -6435.5
$ws.Range("N132").Value = -16134.3638

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 121
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494
# Row 126
$ws.Range("H126").Value = 1458.6666
$ws.Range("I126").Value = 1432.1818
$ws.Range("K126").Value = 4296.5454
$ws.Range("M126").Value = -1826.5454
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
# Row 132
$ws.Range("H132").Value = 3342.3403
$ws.Range("I132").Value = 2691.5945
$ws.Range("J132").Value = 5750.1
$ws.Range("K132").Value = 8074.7835
$ws.Range("L132").Value = 17250.3
$ws.Range("M132").Value = -5544.7835
$ws.Range("N132").Value = -22310.3
# Row 136
$ws.Range("H136").Value = 5409.7
$ws.Range("I136").Value = 1871.1428
$ws.Range("J136").Value = 13666.333
$ws.Range("K136").Value = 5613.428400000001
$ws.Range("L136").Value = 40998.999
$ws.Range("M136").Value = -3063.428400000001

